# Apply the "Add files via upload / Final Presentations with all slides"
# edit: append two new slides (Title and Content layout) at the end of
# the deck, and refresh the cached "datetimeFigureOut" footer field
# (12/3/2019 -> 12/5/2019) on the slide master and every slide layout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Add the two new slides.
#
# The final slide order is:
#   ... existing 4 slides ...
#   5: "Future Plans for PARTY"          (sldId 261)
#   6: "Questions and Where to find us"  (sldId 260)
#
# PowerPoint hands out sldId values in creation order (not list order),
# so to end up with that numbering the "Questions" slide has to be
# created first (claiming the lower id, 260) and then the "Future
# Plans" slide is inserted in front of it (claiming the higher id,
# 261, while still landing on position 5).
# ---------------------------------------------------------------------

$master = $p.SlideMaster
$titleAndContent = $master.CustomLayouts.Item(2)   # "Title and Content"

$insertIndex = $p.Slides.Count + 1

$sQuestions = $p.Slides.AddSlide($insertIndex, $titleAndContent)
$sQuestions.Shapes.Item(1).TextFrame.TextRange.Text = "Questions and Where to find us"
$trQ = $sQuestions.Shapes.Item(2).TextFrame.TextRange
$trQ.Text = "GitHub: "
$trQ.InsertAfter("https://github.com/TSLogan-UTM/PARTYapp")

$sFuture = $p.Slides.AddSlide($insertIndex, $titleAndContent)
$sFuture.Shapes.Item(1).TextFrame.TextRange.Text = "Future Plans for PARTY"
$sFuture.Shapes.Item(2).TextFrame.TextRange.Text = "Better Save & Load Feature`rAdvance Tools working Functionally`rSomewhat better Color control in Advance Mode."

# ---------------------------------------------------------------------
# 2. Refresh the cached date field (Insert > Header & Footer > Apply to
#    All) from 12/3/2019 to 12/5/2019 across the slide master and every
#    slide layout.
# ---------------------------------------------------------------------

function Update-DatePlaceholder($shapes) {
  for ($i = 1; $i -le $shapes.Count; $i++) {
    $shape = $shapes.Item($i)
    if ($shape.HasTextFrame) {
      if ($shape.PlaceholderFormat.Type -eq 16) {
        $shape.TextFrame.TextRange.Text = "12/5/2019"
      }
    }
  }
}

Update-DatePlaceholder $master.Shapes
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
  $layout = $master.CustomLayouts.Item($li)
  Update-DatePlaceholder $layout.Shapes
}
